$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6886.875
$ws.Range("J51").Value = 6942.143
$ws.Range("L51").Value = 6942.143
$ws.Range("N51").Value = -7910.143

$ws.Range("H93").Value = 36000
$ws.Range("J93").Value = 36000
$ws.Range("L93").Value = 36000
$ws.Range("N93").Value = -40992

$ws.Range("H107").Value = 35720080
$ws.Range("I107").Value = 27782622
$ws.Range("K107").Value = 27782622
$ws.Range("M107").Value = -27780702

$ws.Range("H138").Value = 3490.45
$ws.Range("I138").Value = 2466
$ws.Range("J138").Value = 3850.3918
$ws.Range("K138").Value = 7398
$ws.Range("L138").Value = 11551.1754
$ws.Range("M138").Value = -2258
$ws.Range("N138").Value = -21831.1754

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9548.200000000001
$ws.Range("I2").Value = 11364.182
$ws.Range("J2").Value = 7328.6665
$ws.Range("K2").Value = 11364.182
$ws.Range("L2").Value = 7328.6665
$ws.Range("M2").Value = -11251.182
$ws.Range("N2").Value = -7554.6665

$ws.Range("H45").Value = 2060.5557
$ws.Range("I45").Value = 1115.6666
$ws.Range("K45").Value = 1115.6666
$ws.Range("M45").Value = -738.6666

$ws.Range("H74").Value = 4958.403
$ws.Range("I74").Value = 5005.161
$ws.Range("J74").Value = 4378.6
$ws.Range("K74").Value = 5005.161
$ws.Range("L74").Value = 4378.6
$ws.Range("M74").Value = -4131.161
$ws.Range("N74").Value = -6126.6

$ws.Range("H77").Value = 4958.403
$ws.Range("I77").Value = 5005.161
$ws.Range("J77").Value = 4378.6
$ws.Range("K77").Value = 25025.805
$ws.Range("L77").Value = 21893
$ws.Range("M77").Value = -20657.805
$ws.Range("N77").Value = -30629

$ws.Range("H88").Value = 1980.15
$ws.Range("I88").Value = 1783.125
$ws.Range("K88").Value = 1783.125
$ws.Range("M88").Value = -1377.125

$ws.Range("H91").Value = 1980.15
$ws.Range("I91").Value = 1783.125
$ws.Range("K91").Value = 1783.125
$ws.Range("M91").Value = -379.125

$ws.Range("H110").Value = 16488.533
$ws.Range("I110").Value = 17130.6
$ws.Range("J110").Value = 15204.4
$ws.Range("K110").Value = 17130.6
$ws.Range("L110").Value = 15204.4
$ws.Range("M110").Value = -15085.6
$ws.Range("N110").Value = -19294.4

$ws.Range("H116").Value = 9548.200000000001
$ws.Range("I116").Value = 11364.182
$ws.Range("J116").Value = 7328.6665
$ws.Range("K116").Value = 11364.182
$ws.Range("L116").Value = 7328.6665
$ws.Range("M116").Value = -9070.182000000001
$ws.Range("N116").Value = -11916.6665

$ws.Range("H122").Value = 9210.825000000001
$ws.Range("I122").Value = 6658.3105
$ws.Range("K122").Value = 19974.9315
$ws.Range("M122").Value = -17524.9315

$ws.Range("H132").Value = 2806.1592
$ws.Range("I132").Value = 2405.1025
$ws.Range("J132").Value = 5934.4
$ws.Range("K132").Value = 7215.3075
$ws.Range("L132").Value = 17803.2
$ws.Range("M132").Value = -4685.3075
$ws.Range("N132").Value = -22863.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9548.200000000001
$ws.Range("I3").Value = 11364.182
$ws.Range("J3").Value = 7328.6665
$ws.Range("K3").Value = 11364.182
$ws.Range("L3").Value = 7328.6665
$ws.Range("M3").Value = -11250.182
$ws.Range("N3").Value = -7556.6665

$ws.Range("H99").Value = 7315.56
$ws.Range("I99").Value = 10851.417
$ws.Range("K99").Value = 10851.417
$ws.Range("M99").Value = -9353.416999999999

$ws.Range("H105").Value = 2094.077
$ws.Range("I105").Value = 2094.077
$ws.Range("K105").Value = 2094.077
$ws.Range("M105").Value = -347.0770000000002

$ws.Range("H134").Value = 2642.541
$ws.Range("I134").Value = 2914.0476
$ws.Range("K134").Value = 8742.1428
$ws.Range("M134").Value = -6207.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 649.9048
$ws.Range("I22").Value = 486.78946
$ws.Range("J22").Value = 2199.5
$ws.Range("K22").Value = 486.78946
$ws.Range("L22").Value = 2199.5
$ws.Range("M22").Value = -136.78946
$ws.Range("N22").Value = -2899.5

$ws.Range("H31").Value = 53068092
$ws.Range("J31").Value = 90981176
$ws.Range("L31").Value = 90981176
$ws.Range("N31").Value = -90981766

$ws.Range("H34").Value = 53068092
$ws.Range("J34").Value = 90981176
$ws.Range("L34").Value = 90981176
$ws.Range("N34").Value = -90981580

$ws.Range("H41").Value = 104.71429
$ws.Range("I41").Value = 104.71429
$ws.Range("K41").Value = 104.71429
$ws.Range("M41").Value = 323.28571

$ws.Range("H122").Value = 53119
$ws.Range("I122").Value = 3077.9285
$ws.Range("J122").Value = 169881.5
$ws.Range("K122").Value = 9233.7855
$ws.Range("L122").Value = 509644.5
$ws.Range("M122").Value = -6783.7855
$ws.Range("N122").Value = -514544.5

$ws.Range("H132").Value = 3962.5217
$ws.Range("I132").Value = 3530.4285
$ws.Range("K132").Value = 10591.2855
$ws.Range("M132").Value = -8061.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4758.8975
$ws.Range("J132").Value = 2341.0715
$ws.Range("L132").Value = 21069.6435
$ws.Range("N132").Value = -26129.6435

$ws.Range("H137").Value = 37039892
$ws.Range("I137").Value = 2614.5
$ws.Range("J137").Value = 111114450
$ws.Range("K137").Value = 7843.5
$ws.Range("L137").Value = 333343350
$ws.Range("M137").Value = -2743.5
$ws.Range("N137").Value = -333353550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 9000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H70").Value = 125006930
$ws.Range("I70").Value = 5466
$ws.Range("J70").Value = 200007790
$ws.Range("K70").Value = 5466
$ws.Range("L70").Value = 200007790
$ws.Range("M70").Value = -5196
$ws.Range("N70").Value = -200008330

$ws.Range("H73").Value = 125006930
$ws.Range("I73").Value = 5466
$ws.Range("J73").Value = 200007790
$ws.Range("K73").Value = 5466
$ws.Range("L73").Value = 200007790
$ws.Range("M73").Value = -4530
$ws.Range("N73").Value = -200009662

$ws.Range("H97").Value = 2100.25
$ws.Range("J97").Value = 3448.5
$ws.Range("L97").Value = 3448.5
$ws.Range("N97").Value = -4440.5

$ws.Range("H126").Value = 11135.516
$ws.Range("I126").Value = 11559.6875
$ws.Range("K126").Value = 34679.0625
$ws.Range("M126").Value = -32209.0625

$ws.Range("H132").Value = 41298.96
$ws.Range("I132").Value = 46055.914
$ws.Range("J132").Value = 4829
$ws.Range("K132").Value = 138167.742
$ws.Range("L132").Value = 14487
$ws.Range("M132").Value = -135637.742
$ws.Range("N132").Value = -19547

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 289999.28
$ws.Range("J20").Value = 5999.4
$ws.Range("L20").Value = 5999.4
$ws.Range("N20").Value = -6451.4

$ws.Range("H22").Value = 2223.8
$ws.Range("I22").Value = 689
$ws.Range("J22").Value = 3247
$ws.Range("K22").Value = 689
$ws.Range("L22").Value = 3247
$ws.Range("M22").Value = -394
$ws.Range("N22").Value = -3837

$ws.Range("H27").Value = 2223.8
$ws.Range("I27").Value = 689
$ws.Range("J27").Value = 3247
$ws.Range("K27").Value = 689
$ws.Range("L27").Value = 3247
$ws.Range("M27").Value = -582
$ws.Range("N27").Value = -3461

$ws.Range("H132").Value = 6520.913
$ws.Range("I132").Value = 2449.5
$ws.Range("K132").Value = 7348.5
$ws.Range("M132").Value = -4818.5

$ws.Range("H134").Value = 59427
$ws.Range("J134").Value = 59427
$ws.Range("L134").Value = 59427
$ws.Range("N134").Value = -69567

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6199.4
$ws.Range("J81").Value = 7999
$ws.Range("L81").Value = 15998
$ws.Range("N81").Value = -18120

$ws.Range("H84").Value = 6199.4
$ws.Range("J84").Value = 7999
$ws.Range("L84").Value = 79990
$ws.Range("N84").Value = -90598

$ws.Range("H96").Value = 96349.45
$ws.Range("I96").Value = 204070
$ws.Range("K96").Value = 204070
$ws.Range("M96").Value = -202697

$ws.Range("H132").Value = 20410490
$ws.Range("I132").Value = 40002060
$ws.Range("K132").Value = 120006180
$ws.Range("M132").Value = -120003650

$ws.Range("H133").Value = 99499
$ws.Range("J133").Value = 99499
$ws.Range("L133").Value = 99499
$ws.Range("N133").Value = -109619

$ws.Range("H135").Value = 59998.332
$ws.Range("J135").Value = 59998.332
$ws.Range("L135").Value = 59998.332
$ws.Range("N135").Value = -70138.33199999999
